$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- Paragraph "Тринкеты": wrap run in proofErr spellStart/spellEnd ---
$p = $d.Paragraphs(1)
$xml = "<w:p $ns>" +
       "<w:pPr><w:pStyle w:val='a3'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:t>Тринкеты</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# --- Paragraph "Квест на серебряный меч": split "Квест" into its own
#     spell-checked run, keep " на серебряный меч" as a second run ---
$p = $d.Paragraphs(3)
$xml = "<w:p $ns>" +
       "<w:pPr><w:pStyle w:val='a3'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:t>Квест</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> на серебряный меч</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# --- Paragraph "Ачивки": wrap run in proofErr spellStart/spellEnd ---
$p = $d.Paragraphs(7)
$xml = "<w:p $ns>" +
       "<w:pPr><w:pStyle w:val='a3'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:t>Ачивки</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# --- Paragraph "Звуки": wrap run (keeps its en-US rPr) in proofErr ---
$p = $d.Paragraphs(8)
$xml = "<w:p $ns>" +
       "<w:pPr><w:pStyle w:val='a3'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:rPr><w:lang w:val='en-US'/></w:rPr><w:t>Звуки</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# --- Paragraph "Убрать читы и кнопку «Убрать достижения»": re-split runs,
#     close the quote, drop the stray trailing bookmark-only run, and
#     split off a new final "Вернуть ачивки" list item carrying the
#     _GoBack bookmark ---
$p = $d.Paragraphs(12)
$xml = "<w:p $ns>" +
       "<w:pPr><w:pStyle w:val='a3'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>Убрать </w:t></w:r>" +
       "<w:proofErr w:type='spellStart'/>" +
       "<w:r><w:t>читы</w:t></w:r>" +
       "<w:proofErr w:type='spellEnd'/>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:r><w:t>и кнопку &#171;Убрать достижения&#187;</w:t></w:r>" +
       "</w:p>" +
       "<w:p $ns>" +
       "<w:pPr><w:pStyle w:val='a3'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='1'/></w:numPr></w:pPr>" +
       "<w:r><w:t>Вернуть ачивки</w:t></w:r>" +
       "<w:bookmarkStart w:id='0' w:name='_GoBack'/><w:bookmarkEnd w:id='0'/>" +
       "</w:p>"
$p.Range.InsertXML($xml)
